# Update Step3_DataPts_* sheets with new First_Noticeable_Increase_Index (C),
# First_Noticeable_Increase_Cumulative_Value (E), and recomputed Pulse_Width (G)
# values resulting from the new configurable zero_before_threshold parameter
# used to determine the First Rise Point.

$wb = $excel.ActiveWorkbook

# New First_Noticeable_Increase_Index (column C) values per signal segment row (2-6).
# These are the same across all four Step3_DataPts_* sheets.
$newC = @{ 2 = 87; 3 = 89; 4 = 33; 5 = 30; 6 = 41 }

# New First_Noticeable_Increase_Cumulative_Value (column E) values - only the
# 0.5 and 0.7 threshold sheets show a visible change (E is independent of the
# threshold itself, but was only re-written in those two sheets).
$newE = @{ 2 = 0.007460579820718383; 3 = 0.01037312888609691; 4 = 0.0169815007156006; 5 = 0.005888959148646303; 6 = 0.06351611211930204 }

# New Pulse_Width (column G) values per sheet (Point_Exceeds_Index (D) minus
# the new First_Noticeable_Increase_Index (C)).
$newG = @{
    "Step3_DataPts_0.5" = @{ 2 = 48; 3 = 44; 4 = 64; 5 = 62; 6 = 53 }
    "Step3_DataPts_0.7" = @{ 2 = 69; 3 = 67; 4 = 73; 5 = 70; 6 = 66 }
    "Step3_DataPts_0.8" = @{ 2 = 75; 3 = 69; 4 = 82; 5 = 79; 6 = 77 }
    "Step3_DataPts_0.9" = @{ 2 = 85; 3 = 77; 4 = 114; 5 = 114; 6 = 101 }
}

$sheetsWithE = @("Step3_DataPts_0.5", "Step3_DataPts_0.7")

foreach ($sheetName in @("Step3_DataPts_0.5", "Step3_DataPts_0.7", "Step3_DataPts_0.8", "Step3_DataPts_0.9")) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in 2..6) {
        $ws.Cells.Item($row, 3).Value = $newC[$row]              # Column C
        if ($sheetsWithE -contains $sheetName) {
            $ws.Cells.Item($row, 5).Value = $newE[$row]          # Column E
        }
        $ws.Cells.Item($row, 7).Value = $newG[$sheetName][$row]  # Column G
    }
}
